# Append three new "Australian test users" rows (Bank transfer / Voucher /
# TopUp Mobile refactor test data) below the existing data on the Users sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @{ B = "BETA"; C = "70000017358"; D = "uxnvH+test@ogJ.com"; E = "1234567"; F = "Australia" },
    @{ B = "BETA"; C = "70000019235"; D = "qUImP+test@Zcb.com"; E = "1234567"; F = "Australia" },
    @{ B = "PROD"; C = "70000027021"; D = "jtAgW+test@nyb.com"; E = "1234567"; F = "Australia" }
)

$startRow = 18
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    # C (user id) and E (phone number) look like numbers, so force text
    # formatting first -- mirrors how the existing rows store these as
    # shared strings rather than numeric cells.
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 5).NumberFormat = "@"

    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
}
